$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.927.25"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.304.69"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.614"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.293.39"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("E10").Value = "  +5.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.626"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "3.841.70"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.119"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "17.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "3.295.82"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "63.983.45"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.60"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.977"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.61%  "
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.57%  "
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "60.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.00%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "563.38"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.363"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "0.0₃0721"
$ws.Range("E41").Value = "  -3.28%  "
$ws.Range("D42").Value = "3.032.47"
$ws.Range("E42").Value = "  -2.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0410"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.132"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.19%  "
$ws.Range("E50").Value = "  -2.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.79%  "
